# SSDM-12286 Fixed letter case inconsistencies.
#
# "Generated Code Prefix" -> "Generated code prefix"
# "Vocabulary Code"       -> "Vocabulary code"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header label used in the "Property" table (columns shared by both
# SAMPLE_TYPE blocks, rows 12 and 20).
$ws.Range("H12").Value = "Vocabulary code"
$ws.Range("H20").Value = "Vocabulary code"

# Label used in the "SAMPLE_TYPE" summary rows (rows 10 and 18).
$ws.Range("E10").Value = "Generated code prefix"
$ws.Range("E18").Value = "Generated code prefix"

# The saved file had its active selection on E18 (second SAMPLE_TYPE block).
$ws.Range("E18").Select()
